# Updated data to reflect new requirement separation
# - Split the old "Terms Typically Offered" (column D) layout into new
#   Corequisites / Concurrent / Recommended / Terms Typically Offered columns
#   (D, E, F, G respectively), shifting the old "Terms Typically Offered"
#   values into the new column G.
# - A couple of prerequisite descriptions were tweaked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 27

# 1) Move the existing "Terms Typically Offered" values (currently in column D,
#    rows 2-27) over to the new column G before we overwrite column D.
for ($r = 2; $r -le $lastRow; $r++) {
    $oldTerms = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 7).Value2 = $oldTerms
}

# 2) Fill the new Corequisites / Concurrent / Recommended columns (D, E, F)
#    with "NA" for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value2 = "NA"
    $ws.Cells.Item($r, 5).Value2 = "NA"
    $ws.Cells.Item($r, 6).Value2 = "NA"
}

# 3) Update the header row for the new/shifted columns.
$ws.Cells.Item(1, 4).Value2 = "Corequisites"
$ws.Cells.Item(1, 5).Value2 = "Concurrent"
$ws.Cells.Item(1, 6).Value2 = "Recommended"
$ws.Cells.Item(1, 7).Value2 = "Terms Typically Offered"

# 4) Minor text corrections to a few prerequisite descriptions.
$ws.Cells.Item(13, 3).Value2 = "DANC 134 or intermediate level experience as determined by instructor at first class meeting."
$ws.Cells.Item(17, 3).Value2 = "DANC 231 or intermediate level experience as determined by instructor at first class meeting."
$ws.Cells.Item(22, 3).Value2 = "One of the DANC 332, DANC 331, or DANC 345."
